# Fix bugs in data download: rename labels in column A of the TFEC sheet
# (dashboard/Data/EnergyEfficiency.xlsx)
#   A1: "VISUALIZATION"                      -> "Type"
#   A3: "Efficient electricity appliances"    -> "Energy Efficiency residential"
#   A4: "Energy Efficiency in commercial"     -> "Energy Efficiency commercial"
# A2, B1:B4 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Type"
$ws.Range("A3").Value = "Energy Efficiency residential"
$ws.Range("A4").Value = "Energy Efficiency commercial"
